$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: was mon1 / mon1_desc / Mob / ... / H4=mon1(text) / J4=FALSE
# becomes 할아버지 / (cleared) / (cleared) / H4=0 / I4=0 / J4=TRUE / K4=할아버지
$ws.Range("A4").Value = "할아버지"
$ws.Range("B4").ClearContents()
$ws.Range("C4").ClearContents()
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = $true
$ws.Range("K4").Value = "할아버지"

# Row 5: was mon2 / mon2_desc / Mob / ... / H5=mon2(text) / J5=FALSE
# becomes 잼민이 / (cleared) / (cleared) / H5=0 / I5=0 / J5=FALSE
$ws.Range("A5").Value = "잼민이"
$ws.Range("B5").ClearContents()
$ws.Range("C5").ClearContents()
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = $false

# Rows 6 and 7 are removed entirely (revert of merged monster rows)
$ws.Range("A6:K7").EntireRow.Delete()

# Sheet view / dimension bookkeeping
$ws.Range("J10").Select()
